$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8683426976203918
$ws.Range("B1").Value = 2.622020483016968
$ws.Range("C1").Value = 4.741023063659668
$ws.Range("D1").Value = 2.221381902694702
$ws.Range("E1").Value = 1.311357736587524
